$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$row = 44
$ws.Cells.Item($row, 1).Value = "Height"
$ws.Cells.Item($row, 2).Value = "float"
$ws.Cells.Item($row, 3).Value = $false
$ws.Cells.Item($row, 4).Value = $false
$ws.Cells.Item($row, 5).Value = $false
$ws.Cells.Item($row, 6).Value = $true
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = "Friend"
$ws.Cells.Item($row, 10).Value = "模型高度"
$ws.Cells.Item($row, 9).NumberFormat = "@"
